$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3453645
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 3565043.2
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 10695129.6
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -10695465.6

$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H40").Value = 2050
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 2111.111
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2111.111
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2461.111

$ws.Range("H52").Value = 1000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H58").Value = 443.18182
$ws.Range("I58").Value = 206.44444
$ws.Range("J58").Value = 1508.5
$ws.Range("K58").Value = 619.33332
$ws.Range("L58").Value = 4525.5
$ws.Range("M58").Value = -469.33332
$ws.Range("N58").Value = -4825.5

$ws.Range("H106").Value = 26788.75
$ws.Range("I106").Value = 26788.75
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 26788.75
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H112").Value = 1873.8276
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 1926.4642
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 5779.392599999999
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -7995.392599999999

$ws.Range("H113").Value = 2062.0454
$ws.Range("I113").Value = 1805
$ws.Range("J113").Value = 2612.8572
$ws.Range("K113").Value = 1805
$ws.Range("L113").Value = 2612.8572
$ws.Range("M113").Value = 1449
$ws.Range("N113").Value = -9120.8572

$ws.Range("H116").Value = 2000.909
$ws.Range("I116").Value = 2063.75
$ws.Range("J116").Value = 1833.3334
$ws.Range("K116").Value = 2063.75
$ws.Range("L116").Value = 1833.3334
$ws.Range("M116").Value = 1378.25
$ws.Range("N116").Value = -8717.3334

$ws.Range("H127").Value = 719016.3
$ws.Range("I127").Value = 365
$ws.Range("J127").Value = 1198117.2
$ws.Range("K127").Value = 1095
$ws.Range("L127").Value = 3594351.6
$ws.Range("M127").Value = 3865
$ws.Range("N127").Value = -3604271.6

$ws.Range("H129").Value = 828.42426
$ws.Range("I129").Value = 485.27274
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1455.81822
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3544.18178
$ws.Range("N129").Value = -13000

$ws.Range("H132").Value = 3970340.8
$ws.Range("I132").Value = 1557.5636
$ws.Range("J132").Value = 31255724
$ws.Range("K132").Value = 4672.6908
$ws.Range("L132").Value = 93767172
$ws.Range("M132").Value = -2142.6908
$ws.Range("N132").Value = -93772232

$ws.Range("H138").Value = 12347784
$ws.Range("J138").Value = 4147.7144
$ws.Range("L138").Value = 12443.1432
$ws.Range("N138").Value = -22723.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1447.0625
$ws.Range("I2").Value = 1085.4546
$ws.Range("K2").Value = 1085.4546
$ws.Range("M2").Value = -972.4546

$ws.Range("H32").Value = 9600.67
$ws.Range("I32").Value = 9073.272000000001
$ws.Range("J32").Value = 11631.15
$ws.Range("K32").Value = 9073.272000000001
$ws.Range("L32").Value = 11631.15
$ws.Range("M32").Value = -8786.272000000001
$ws.Range("N32").Value = -12205.15

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H116").Value = 1447.0625
$ws.Range("I116").Value = 1085.4546
$ws.Range("K116").Value = 1085.4546
$ws.Range("M116").Value = 1208.5454

$ws.Range("H117").Value = 23557
$ws.Range("J117").Value = 23557
$ws.Range("L117").Value = 23557
$ws.Range("N117").Value = -32735

$ws.Range("H122").Value = 16336
$ws.Range("I122").Value = 17878.25
$ws.Range("K122").Value = 53634.75
$ws.Range("M122").Value = -51184.75

$ws.Range("H134").Value = 130000
$ws.Range("J134").Value = 130000
$ws.Range("L134").Value = 130000
$ws.Range("N134").Value = -140140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1447.0625
$ws.Range("I3").Value = 1085.4546
$ws.Range("K3").Value = 1085.4546
$ws.Range("M3").Value = -971.4546

$ws.Range("H68").Value = 42047.5
$ws.Range("J68").Value = 42047.5
$ws.Range("L68").Value = 42047.5
$ws.Range("N68").Value = -43669.5

$ws.Range("H71").Value = 42047.5
$ws.Range("J71").Value = 42047.5
$ws.Range("L71").Value = 126142.5
$ws.Range("N71").Value = -134254.5

$ws.Range("H94").Value = 862.95
$ws.Range("I94").Value = 767.0714
$ws.Range("J94").Value = 1086.6666
$ws.Range("K94").Value = 767.0714
$ws.Range("L94").Value = 1086.6666
$ws.Range("M94").Value = -316.0714
$ws.Range("N94").Value = -1988.6666

$ws.Range("H134").Value = 3015.257
$ws.Range("I134").Value = 2472.4614
$ws.Range("J134").Value = 4583.3335
$ws.Range("K134").Value = 7417.3842
$ws.Range("L134").Value = 13750.0005
$ws.Range("M134").Value = -4882.3842
$ws.Range("N134").Value = -18820.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8055.4
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 9069.25
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 9069.25
$ws.Range("M51").Value = -3264
$ws.Range("N51").Value = -10541.25

$ws.Range("H61").Value = 8055.4
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 9069.25
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 9069.25
$ws.Range("M61").Value = -3652
$ws.Range("N61").Value = -9765.25

$ws.Range("H132").Value = 14287173
$ws.Range("I132").Value = 16667894
$ws.Range("J132").Value = 2846
$ws.Range("K132").Value = 50003682
$ws.Range("L132").Value = 8538
$ws.Range("M132").Value = -50001152
$ws.Range("N132").Value = -13598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1425
$ws.Range("I5").Value = 282.14285
$ws.Range("J5").Value = 4425
$ws.Range("K5").Value = 846.4285500000001
$ws.Range("L5").Value = 13275
$ws.Range("M5").Value = -734.4285500000001
$ws.Range("N5").Value = -13499

$ws.Range("H122").Value = 661.7619
$ws.Range("I122").Value = 488.2353
$ws.Range("J122").Value = 1399.25
$ws.Range("K122").Value = 4394.1177
$ws.Range("L122").Value = 12593.25
$ws.Range("M122").Value = -1944.1177
$ws.Range("N122").Value = -17493.25

$ws.Range("H135").Value = 1425
$ws.Range("I135").Value = 282.14285
$ws.Range("J135").Value = 4425
$ws.Range("K135").Value = 2539.28565
$ws.Range("L135").Value = 39825
$ws.Range("M135").Value = -4.28565000000026
$ws.Range("N135").Value = -44895

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12555.571
$ws.Range("J57").Value = 23996.666
$ws.Range("L57").Value = 23996.666
$ws.Range("N57").Value = -25636.666

$ws.Range("H58").Value = 13620
$ws.Range("I58").Value = 9433.333000000001
$ws.Range("J58").Value = 19900
$ws.Range("K58").Value = 9433.333000000001
$ws.Range("L58").Value = 19900
$ws.Range("M58").Value = -9156.333000000001
$ws.Range("N58").Value = -20454

$ws.Range("H118").Value = 14096.774
$ws.Range("J118").Value = 14096.774
$ws.Range("L118").Value = 14096.774
$ws.Range("N118").Value = -17410.774

$ws.Range("H132").Value = 5943.0645
$ws.Range("I132").Value = 6425.1665
$ws.Range("J132").Value = 4290.143
$ws.Range("K132").Value = 19275.4995
$ws.Range("L132").Value = 12870.429
$ws.Range("M132").Value = -16745.4995
$ws.Range("N132").Value = -17930.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H58").Value = 2650
$ws.Range("I58").Value = 200
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 200
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = 60
$ws.Range("N58").Value = -10520

$ws.Range("H132").Value = 7939721
$ws.Range("I132").Value = 2734.6667
$ws.Range("J132").Value = 23813692
$ws.Range("K132").Value = 8204.000100000001
$ws.Range("L132").Value = 71441076
$ws.Range("M132").Value = -5674.000100000001
$ws.Range("N132").Value = -71446136

$ws.Range("H136").Value = 41668892
$ws.Range("I136").Value = 50002270
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 150006810
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -150004260
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7319.6
$ws.Range("I62").Value = 6900
$ws.Range("J62").Value = 7599.3335
$ws.Range("K62").Value = 6900
$ws.Range("L62").Value = 7599.3335
$ws.Range("M62").Value = -6276
$ws.Range("N62").Value = -8847.333500000001

$ws.Range("H65").Value = 7319.6
$ws.Range("I65").Value = 6900
$ws.Range("J65").Value = 7599.3335
$ws.Range("K65").Value = 34500
$ws.Range("L65").Value = 37996.6675
$ws.Range("M65").Value = -31380
$ws.Range("N65").Value = -44236.6675

$ws.Range("H113").Value = 1483.4166
$ws.Range("I113").Value = 200.5
$ws.Range("J113").Value = 2124.875
$ws.Range("K113").Value = 601.5
$ws.Range("L113").Value = 6374.625
$ws.Range("M113").Value = 1568.5
$ws.Range("N113").Value = -10714.625

$ws.Range("H132").Value = 1537.8387
$ws.Range("I132").Value = 1427.8572
$ws.Range("J132").Value = 1768.8
$ws.Range("K132").Value = 4283.571599999999
$ws.Range("L132").Value = 5306.4
$ws.Range("M132").Value = -1753.571599999999
$ws.Range("N132").Value = -10366.4

$ws.Range("H136").Value = 1630.591
$ws.Range("I136").Value = 1478.8667
$ws.Range("K136").Value = 4436.6001
$ws.Range("M136").Value = -1886.6001

